# Edit script: insert a new weekly record for "Orégano" at row 27,
# pushing the existing rows 27-31 down to 28-32 (their values stay the same).
# The new row 27 gets the fresh data point (date 44468 / 2021-09-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 27; this shifts rows 27:31 -> 28:32,
# carrying their existing cell values/styles along.
$ws.Rows.Item(27).Insert()

# The freshly inserted row 27 is blank; copy the formatting of the date cell
# right below it (the row that used to be 27, now 28) so that the new date
# cell (D27) keeps the same date/time number format.
$ws.Cells.Item(28, 4).Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 27 with the new data point.
$ws.Cells.Item(27, 1).Value = 9
$ws.Cells.Item(27, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(27, 3).Value = "Metropolitana"
$ws.Cells.Item(27, 4).Value = 44468
$ws.Cells.Item(27, 5).Value = 13
$ws.Cells.Item(27, 6).Value = 100112029
$ws.Cells.Item(27, 7).Value = "Orégano"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 16
$ws.Cells.Item(27, 11).Value = 10000
$ws.Cells.Item(27, 12).Value = 11000
$ws.Cells.Item(27, 13).Value = 10500
$ws.Cells.Item(27, 14).Value = "$/docena de atados"
$ws.Cells.Item(27, 15).Value = "Región Metropolitana"
$ws.Cells.Item(27, 16).Value = 3500
$ws.Cells.Item(27, 17).Value = 3
$ws.Cells.Item(27, 18).Value = "Hortaliza"
